$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.202.69"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "1.969.84"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'248.12"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D7").Value = "'0.4891"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").Value = "'44.53"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'0.2965"
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").Value = "'0.06854"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").Value = "'19.26"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "'106.88"
$ws.Range("E12").Value = "  -3.83%  "
$ws.Range("D13").Value = "1.958.64"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").Value = "'0.07782"
$ws.Range("E14").Value = "  +2.89%  "
$ws.Range("D15").Value = "'5.451"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "'0.7164"
$ws.Range("E16").Value = "  +7.04%  "
$ws.Range("D17").Value = "'285.40"
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("D18").Value = "31.079.35"
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("D19").Value = "'13.35"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000007773"
$ws.Range("E20").Value = "  +2.53%  "
$ws.Range("B21").Value = "BitDAO"
$ws.Range("C21").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D21").Value = "'0.4957"
$ws.Range("E21").Value = "  +12.19%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.648"
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.218.31"
$ws.Range("E23").Value = "  +2.66%  "
$ws.Range("D24").Value = "'0.9977"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "'0.9996"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "'6.651"
$ws.Range("E26").Value = "  +3.66%  "
$ws.Range("D27").Value = "'10.07"
$ws.Range("E27").Value = "  +6.42%  "
$ws.Range("D28").Value = "'170.09"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("D29").Value = "'20.10"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "'2.208"
$ws.Range("E30").Value = "  +6.19%  "
$ws.Range("D31").Value = "'0.1071"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").Value = "'1.451"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("D33").Value = "'4.805"
$ws.Range("E33").Value = "  +18.75%  "
$ws.Range("D34").Value = "'4.530"
$ws.Range("E34").Value = "  +9.53%  "
$ws.Range("D35").Value = "'0.05084"
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("D36").Value = "'0.7724"
$ws.Range("E36").Value = "  +4.94%  "
$ws.Range("D37").Value = "'1.174"
$ws.Range("E37").Value = "  +3.58%  "
$ws.Range("D38").Value = "'0.02056"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("D39").Value = "'2.739"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").Value = "'2.720"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").Value = "'2.134"
$ws.Range("E41").Value = "  +5.67%  "
$ws.Range("D42").Value = "'6.435"
$ws.Range("E42").Value = "  +11.62%  "
$ws.Range("D43").Value = "'74.12"
$ws.Range("E43").Value = "  +7.11%  "
$ws.Range("D44").Value = "'0.8886"
$ws.Range("E44").Value = "  +2.78%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4488"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'109.96"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.000.17"
$ws.Range("E48").Value = "  +18.91%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.508"
$ws.Range("E49").Value = "  +4.16%  "
$ws.Range("D50").Value = "'0.1272"
$ws.Range("E50").Value = "  +3.60%  "
$ws.Range("D51").Value = "'9.427"
$ws.Range("E51").Value = "  +2.71%  "
